# =====================================================================
# feat: add 2022-Q4 data
#
# - Insert a new worksheet "2022-Q4" immediately after "总计", holding
#   the per-fund holding detail for the new quarter (same B..H schema
#   as the existing quarterly sheets).
# - Update the "总计" (summary) sheet: add a new top data row for
#   2022-Q4 and shift the existing quarterly summary rows down by one,
#   keeping the sequential 0-based index in column A.
# =====================================================================

$wb = $excel.ActiveWorkbook
$summarySheet = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------
# 1. Create the new "2022-Q4" worksheet right after "总计"; all the
#    other quarterly sheets keep their own name/content and simply
#    shift one tab to the right.
# ---------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Add($null, $summarySheet)
$wsQ4.Name = "2022-Q4"
# Match the outline defaults ("summary rows below / summary columns to the
# right") that every other sheet in this workbook already carries.
$wsQ4.Outline.SummaryRow = 1
$wsQ4.Outline.SummaryColumn = 1

# Header row (B1:H1)
$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"

# Data rows 2..20: column A = 0-based row index (number), H = rank
# (number). Columns B..G are stored as TEXT even though several look
# numeric (fund code leading zeros, "25.00"-style figures) -- this
# matches how the source data is typed, so each is entered with a
# leading apostrophe to force Excel's text (quote-prefix) storage.
$wsQ4.Cells.Item(2, 1).Value = 0
$wsQ4.Cells.Item(2, 2).Value = "'" + "001933"
$wsQ4.Cells.Item(2, 3).Value = "华商新兴活力灵活配置混合"
$wsQ4.Cells.Item(2, 4).Value = "'" + "25.00"
$wsQ4.Cells.Item(2, 5).Value = "'" + "90.39"
$wsQ4.Cells.Item(2, 6).Value = "'" + "3.99"
$wsQ4.Cells.Item(2, 7).Value = "'" + "0.9975"
$wsQ4.Cells.Item(2, 8).Value = 6

$wsQ4.Cells.Item(3, 1).Value = 1
$wsQ4.Cells.Item(3, 2).Value = "'" + "010550"
$wsQ4.Cells.Item(3, 3).Value = "华商双擎领航混合"
$wsQ4.Cells.Item(3, 4).Value = "'" + "12.41"
$wsQ4.Cells.Item(3, 5).Value = "'" + "90.98"
$wsQ4.Cells.Item(3, 6).Value = "'" + "3.99"
$wsQ4.Cells.Item(3, 7).Value = "'" + "0.4952"
$wsQ4.Cells.Item(3, 8).Value = 8

$wsQ4.Cells.Item(4, 1).Value = 2
$wsQ4.Cells.Item(4, 2).Value = "'" + "013886"
$wsQ4.Cells.Item(4, 3).Value = "华商新能源汽车混合A"
$wsQ4.Cells.Item(4, 4).Value = "'" + "9.05"
$wsQ4.Cells.Item(4, 5).Value = "'" + "89.00"
$wsQ4.Cells.Item(4, 6).Value = "'" + "4.00"
$wsQ4.Cells.Item(4, 7).Value = "'" + "0.3620"
$wsQ4.Cells.Item(4, 8).Value = 8

$wsQ4.Cells.Item(5, 1).Value = 3
$wsQ4.Cells.Item(5, 2).Value = "'" + "506001"
$wsQ4.Cells.Item(5, 3).Value = "万家科创板 2 年定期开放混合"
$wsQ4.Cells.Item(5, 4).Value = "'" + "6.18"
$wsQ4.Cells.Item(5, 5).Value = "'" + "95.02"
$wsQ4.Cells.Item(5, 6).Value = "'" + "4.71"
$wsQ4.Cells.Item(5, 7).Value = "'" + "0.2911"
$wsQ4.Cells.Item(5, 8).Value = 4

$wsQ4.Cells.Item(6, 1).Value = 4
$wsQ4.Cells.Item(6, 2).Value = "'" + "005358"
$wsQ4.Cells.Item(6, 3).Value = "东方阿尔法精选灵活配置混合A"
$wsQ4.Cells.Item(6, 4).Value = "'" + "3.45"
$wsQ4.Cells.Item(6, 5).Value = "'" + "94.28"
$wsQ4.Cells.Item(6, 6).Value = "'" + "7.36"
$wsQ4.Cells.Item(6, 7).Value = "'" + "0.2539"
$wsQ4.Cells.Item(6, 8).Value = 7

$wsQ4.Cells.Item(7, 1).Value = 5
$wsQ4.Cells.Item(7, 2).Value = "'" + "010616"
$wsQ4.Cells.Item(7, 3).Value = "国金自主创新混合C"
$wsQ4.Cells.Item(7, 4).Value = "'" + "3.24"
$wsQ4.Cells.Item(7, 5).Value = "'" + "94.49"
$wsQ4.Cells.Item(7, 6).Value = "'" + "5.54"
$wsQ4.Cells.Item(7, 7).Value = "'" + "0.1795"
$wsQ4.Cells.Item(7, 8).Value = 10

$wsQ4.Cells.Item(8, 1).Value = 6
$wsQ4.Cells.Item(8, 2).Value = "'" + "162207"
$wsQ4.Cells.Item(8, 3).Value = "泰达宏利效率优选混合（LOF）"
$wsQ4.Cells.Item(8, 4).Value = "'" + "4.72"
$wsQ4.Cells.Item(8, 5).Value = "'" + "68.20"
$wsQ4.Cells.Item(8, 6).Value = "'" + "3.65"
$wsQ4.Cells.Item(8, 7).Value = "'" + "0.1723"
$wsQ4.Cells.Item(8, 8).Value = 10

$wsQ4.Cells.Item(9, 1).Value = 7
$wsQ4.Cells.Item(9, 2).Value = "'" + "010615"
$wsQ4.Cells.Item(9, 3).Value = "国金自主创新混合A"
$wsQ4.Cells.Item(9, 4).Value = "'" + "3.08"
$wsQ4.Cells.Item(9, 5).Value = "'" + "94.49"
$wsQ4.Cells.Item(9, 6).Value = "'" + "5.54"
$wsQ4.Cells.Item(9, 7).Value = "'" + "0.1706"
$wsQ4.Cells.Item(9, 8).Value = 10

$wsQ4.Cells.Item(10, 1).Value = 8
$wsQ4.Cells.Item(10, 2).Value = "'" + "013887"
$wsQ4.Cells.Item(10, 3).Value = "华商新能源汽车混合C"
$wsQ4.Cells.Item(10, 4).Value = "'" + "3.70"
$wsQ4.Cells.Item(10, 5).Value = "'" + "89.00"
$wsQ4.Cells.Item(10, 6).Value = "'" + "4.00"
$wsQ4.Cells.Item(10, 7).Value = "'" + "0.1480"
$wsQ4.Cells.Item(10, 8).Value = 8

$wsQ4.Cells.Item(11, 1).Value = 9
$wsQ4.Cells.Item(11, 2).Value = "'" + "014185"
$wsQ4.Cells.Item(11, 3).Value = "招商专精特新股票A"
$wsQ4.Cells.Item(11, 4).Value = "'" + "3.30"
$wsQ4.Cells.Item(11, 5).Value = "'" + "87.72"
$wsQ4.Cells.Item(11, 6).Value = "'" + "3.76"
$wsQ4.Cells.Item(11, 7).Value = "'" + "0.1241"
$wsQ4.Cells.Item(11, 8).Value = 10

$wsQ4.Cells.Item(12, 1).Value = 10
$wsQ4.Cells.Item(12, 2).Value = "'" + "014350"
$wsQ4.Cells.Item(12, 3).Value = "华商卓越成长一年持有混合A"
$wsQ4.Cells.Item(12, 4).Value = "'" + "3.05"
$wsQ4.Cells.Item(12, 5).Value = "'" + "93.10"
$wsQ4.Cells.Item(12, 6).Value = "'" + "3.84"
$wsQ4.Cells.Item(12, 7).Value = "'" + "0.1171"
$wsQ4.Cells.Item(12, 8).Value = 8

$wsQ4.Cells.Item(13, 1).Value = 11
$wsQ4.Cells.Item(13, 2).Value = "'" + "014186"
$wsQ4.Cells.Item(13, 3).Value = "招商专精特新股票C"
$wsQ4.Cells.Item(13, 4).Value = "'" + "2.50"
$wsQ4.Cells.Item(13, 5).Value = "'" + "87.72"
$wsQ4.Cells.Item(13, 6).Value = "'" + "3.76"
$wsQ4.Cells.Item(13, 7).Value = "'" + "0.0940"
$wsQ4.Cells.Item(13, 8).Value = 10

$wsQ4.Cells.Item(14, 1).Value = 12
$wsQ4.Cells.Item(14, 2).Value = "'" + "009467"
$wsQ4.Cells.Item(14, 3).Value = "红土创新科技创新3个月定开混合A"
$wsQ4.Cells.Item(14, 4).Value = "'" + "1.62"
$wsQ4.Cells.Item(14, 5).Value = "'" + "88.86"
$wsQ4.Cells.Item(14, 6).Value = "'" + "3.51"
$wsQ4.Cells.Item(14, 7).Value = "'" + "0.0569"
$wsQ4.Cells.Item(14, 8).Value = 10

$wsQ4.Cells.Item(15, 1).Value = 13
$wsQ4.Cells.Item(15, 2).Value = "'" + "013250"
$wsQ4.Cells.Item(15, 3).Value = "红土创新智能制造混合"
$wsQ4.Cells.Item(15, 4).Value = "'" + "1.19"
$wsQ4.Cells.Item(15, 5).Value = "'" + "90.41"
$wsQ4.Cells.Item(15, 6).Value = "'" + "4.76"
$wsQ4.Cells.Item(15, 7).Value = "'" + "0.0566"
$wsQ4.Cells.Item(15, 8).Value = 6

$wsQ4.Cells.Item(16, 1).Value = 14
$wsQ4.Cells.Item(16, 2).Value = "'" + "010375"
$wsQ4.Cells.Item(16, 3).Value = "国金鑫悦经济新动能混合A"
$wsQ4.Cells.Item(16, 4).Value = "'" + "0.89"
$wsQ4.Cells.Item(16, 5).Value = "'" + "92.34"
$wsQ4.Cells.Item(16, 6).Value = "'" + "5.30"
$wsQ4.Cells.Item(16, 7).Value = "'" + "0.0472"
$wsQ4.Cells.Item(16, 8).Value = 10

$wsQ4.Cells.Item(17, 1).Value = 15
$wsQ4.Cells.Item(17, 2).Value = "'" + "005359"
$wsQ4.Cells.Item(17, 3).Value = "东方阿尔法精选灵活配置混合C"
$wsQ4.Cells.Item(17, 4).Value = "'" + "0.57"
$wsQ4.Cells.Item(17, 5).Value = "'" + "94.28"
$wsQ4.Cells.Item(17, 6).Value = "'" + "7.36"
$wsQ4.Cells.Item(17, 7).Value = "'" + "0.0420"
$wsQ4.Cells.Item(17, 8).Value = 7

$wsQ4.Cells.Item(18, 1).Value = 16
$wsQ4.Cells.Item(18, 2).Value = "'" + "010376"
$wsQ4.Cells.Item(18, 3).Value = "国金鑫悦经济新动能混合C"
$wsQ4.Cells.Item(18, 4).Value = "'" + "0.43"
$wsQ4.Cells.Item(18, 5).Value = "'" + "92.34"
$wsQ4.Cells.Item(18, 6).Value = "'" + "5.30"
$wsQ4.Cells.Item(18, 7).Value = "'" + "0.0228"
$wsQ4.Cells.Item(18, 8).Value = 10

$wsQ4.Cells.Item(19, 1).Value = 17
$wsQ4.Cells.Item(19, 2).Value = "'" + "013173"
$wsQ4.Cells.Item(19, 3).Value = "红土创新科技创新3个月定开混合C"
$wsQ4.Cells.Item(19, 4).Value = "'" + "0.43"
$wsQ4.Cells.Item(19, 5).Value = "'" + "88.86"
$wsQ4.Cells.Item(19, 6).Value = "'" + "3.51"
$wsQ4.Cells.Item(19, 7).Value = "'" + "0.0151"
$wsQ4.Cells.Item(19, 8).Value = 10

$wsQ4.Cells.Item(20, 1).Value = 18
$wsQ4.Cells.Item(20, 2).Value = "'" + "014351"
$wsQ4.Cells.Item(20, 3).Value = "华商卓越成长一年持有混合C"
$wsQ4.Cells.Item(20, 4).Value = "'" + "0.10"
$wsQ4.Cells.Item(20, 5).Value = "'" + "93.10"
$wsQ4.Cells.Item(20, 6).Value = "'" + "3.84"
$wsQ4.Cells.Item(20, 7).Value = "'" + "0.0038"
$wsQ4.Cells.Item(20, 8).Value = 8

# Re-apply the same bold/bordered/centered header style, and the same
# index-column style, used by every other quarterly sheet -- copied
# from "总计" so no new style entries are introduced.
$summarySheet.Range("B1").Copy() | Out-Null
$wsQ4.Range("B1:H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$summarySheet.Range("A2").Copy() | Out-Null
$wsQ4.Range("A2:A20").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# 2. Update the "总计" summary sheet: new 2022-Q4 row on top, every
#    other quarter shifts down one row (index column stays 0..5).
# ---------------------------------------------------------------
$summarySheet.Cells.Item(2, 1).Value = 0
$summarySheet.Cells.Item(2, 2).Value = "2022-Q4"
$summarySheet.Cells.Item(2, 3).Value = 19
$summarySheet.Cells.Item(2, 4).Value = 3.65

$summarySheet.Cells.Item(3, 1).Value = 1
$summarySheet.Cells.Item(3, 2).Value = "2022-Q3"
$summarySheet.Cells.Item(3, 3).Value = 45
$summarySheet.Cells.Item(3, 4).Value = 7.55

$summarySheet.Cells.Item(4, 1).Value = 2
$summarySheet.Cells.Item(4, 2).Value = "2022-Q2"
$summarySheet.Cells.Item(4, 3).Value = 4
$summarySheet.Cells.Item(4, 4).Value = 0.88

$summarySheet.Cells.Item(5, 1).Value = 3
$summarySheet.Cells.Item(5, 2).Value = "2022-Q1"
$summarySheet.Cells.Item(5, 3).Value = 2
$summarySheet.Cells.Item(5, 4).Value = 1.04

$summarySheet.Cells.Item(6, 1).Value = 4
$summarySheet.Cells.Item(6, 2).Value = "2021-Q4"
$summarySheet.Cells.Item(6, 3).Value = 4
$summarySheet.Cells.Item(6, 4).Value = 0.96

$summarySheet.Cells.Item(7, 1).Value = 5
$summarySheet.Cells.Item(7, 2).Value = "2021-Q3"
$summarySheet.Cells.Item(7, 3).Value = 5
$summarySheet.Cells.Item(7, 4).Value = 0.98

# Extend the new row 7 index-column style to match the rest of the
# column (copy from an existing styled cell in the same column).
$summarySheet.Range("A6").Copy() | Out-Null
$summarySheet.Range("A7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
